$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the now-removed "Running_shoes" values for rows 4 and 5
$ws.Range("K4").ClearContents()
$ws.Range("K5").ClearContents()

# Delete rows 11 and 12 entirely (shifts dimension back to A1:M10)
$ws.Rows("11:12").Delete()

# Update the active selection
$ws.Range("D13").Select()
